$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SRPbVT")

# Rename existing categories to "passenger ..."
$ws.Range("A2").Value = "passenger LDVs"
$ws.Range("A3").Value = "passenger HDVs"
$ws.Range("A4").Value = "passenger aircraft"
$ws.Range("A5").Value = "passenger rail"
$ws.Range("A6").Value = "passenger ships"
$ws.Range("A7").Value = "passenger motorbikes"

# Add freight categories in rows 8-13, mirroring rows 2-7
$ws.Range("A8").Value = "freight LDVs"
$ws.Range("A9").Value = "freight HDVs"
$ws.Range("A10").Value = "freight aircraft"
$ws.Range("A11").Value = "freight rail"
$ws.Range("A12").Value = "freight ships"
$ws.Range("A13").Value = "freight motorbikes"

$ws.Range("B8:M8").Formula = "=B2"
$ws.Range("B9:M13").Formula = "=B3"
